# Arbeitszeit.xlsx edit: add a work entry for Juni (June), row 7 (03.06.2021)
#   - Kommt 1 (arrival)   D7 = 14:00
#   - Geht 1  (departure) E7 = 20:30
#   - Bemerkungen (remark) O7 = "Coding"
# and move the active-cell selection on that sheet to G8, matching the
# author's final cursor position when the edit was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Juni")
$ws.Activate()

$ws.Range("D7").Value = 0.58333333333333337
$ws.Range("E7").Value = 0.85416666666666663
$ws.Range("O7").Value = "Coding"

[void]$ws.Range("G8").Select()
